# Update "想去人数" (F column) figures across the four sheets of the
# 北京-漫展信息 workbook, matching the data refresh captured in the diff.
#
# NOTE: this COM-interop runtime does not reliably bind PowerShell named
# parameters (e.g. `-SheetName foo`), so helper functions below use
# positional parameters only.

$wb = $excel.ActiveWorkbook

function Set-FValues {
    param(
        $SheetName,
        $Updates
    )

    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($cellRef in $Updates.Keys) {
        $ws.Range($cellRef).Value = $Updates[$cellRef]
    }
}

Set-FValues "展览" @{
    "F2"  = 297
    "F4"  = 38
    "F5"  = 45
    "F7"  = 1239
    "F8"  = 384
    "F10" = 352
    "F11" = 8226
    "F13" = 9917
    "F17" = 8
    "F23" = 28
    "F25" = 12
    "F26" = 390
    "F27" = 1701
    "F28" = 43
    "F29" = 493
    "F31" = 274
    "F34" = 998
    "F35" = 14
    "F36" = 38
    "F38" = 417
    "F39" = 326
    "F40" = 10
    "F44" = 63
    "F45" = 260
    "F46" = 108
    "F48" = 27
    "F49" = 28
}

Set-FValues "演出" @{
    "F5"  = 99
    "F15" = 54
    "F20" = 363
}

Set-FValues "本地生活" @{
    "F3" = 2767
    "F4" = 331
    "F5" = 194
}

Set-FValues "全部类型" @{
    "F2"  = 297
    "F5"  = 331
    "F6"  = 194
    "F7"  = 38
    "F8"  = 45
    "F10" = 1239
    "F11" = 384
    "F15" = 99
    "F16" = 8226
    "F18" = 9917
    "F20" = 8
    "F23" = 28
    "F24" = 1701
    "F25" = 43
    "F27" = 274
    "F33" = 38
    "F37" = 417
    "F38" = 54
    "F39" = 326
    "F43" = 63
    "F46" = 363
    "F48" = 27
    "F49" = 28
}
